# Sync automático del tracker - adds 3 new "Pending" prediction rows (197-199)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 197; Date = "2025-09-23"; Liga = "La Liga"; Local = "Athletic Club"; Visitante = "Girona"; Prediccion = "Home Win"; Probabilidad = "87.31%"; Cuota = 1.45; EV = "25.34%"; Stake = 3.2; StakePct = 0.05; KellyFrac = 0.5912082196903051 },
    @{ Row = 198; Date = "2025-09-23"; Liga = "Jupiler Pro League"; Local = "Anderlecht"; Visitante = "Gent"; Prediccion = "Home Win"; Probabilidad = "69.62%"; Cuota = 1.8; EV = "24.06%"; Stake = 2; StakePct = 0.03163538932465376; KellyFrac = 0.3163538932465376 },
    @{ Row = 199; Date = "2025-09-23"; Liga = "La Liga"; Local = "Sevilla"; Visitante = "Villarreal"; Prediccion = "Away Win"; Probabilidad = "49.72%"; Cuota = 2.25; EV = "10.74%"; Stake = 0.6; StakePct = 0.009488641199056218; KellyFrac = 0.09488641199056218 }
)

foreach ($r in $rows) {
    $i = $r.Row

    $addrA = "A" + $i
    $ws.Range($addrA).NumberFormat = "@"
    $ws.Range($addrA).Value = $r.Date
    $ws.Range($addrA).Style = "Normal"

    $addrB = "B" + $i
    $ws.Range($addrB).Value = $r.Liga

    $addrC = "C" + $i
    $ws.Range($addrC).Value = $r.Local

    $addrD = "D" + $i
    $ws.Range($addrD).Value = $r.Visitante

    $addrE = "E" + $i
    $ws.Range($addrE).Value = $r.Prediccion

    $addrF = "F" + $i
    $ws.Range($addrF).NumberFormat = "@"
    $ws.Range($addrF).Value = $r.Probabilidad
    $ws.Range($addrF).Style = "Normal"

    $addrG = "G" + $i
    $ws.Range($addrG).Value = $r.Cuota

    $addrH = "H" + $i
    $ws.Range($addrH).NumberFormat = "@"
    $ws.Range($addrH).Value = $r.EV
    $ws.Range($addrH).Style = "Normal"

    $addrI = "I" + $i
    $ws.Range($addrI).Value = $r.Stake

    $addrJ = "J" + $i
    $ws.Range($addrJ).Value = $r.StakePct

    $addrK = "K" + $i
    $ws.Range($addrK).Value = $r.KellyFrac

    $addrL = "L" + $i
    $ws.Range($addrL).Value = "Pending"
}
